# The authored change swaps the contents of ppt/theme/theme1.xml (the
# presentation's "Integral" theme, used by the slide master / all slides)
# and ppt/theme/theme2.xml (an "Office Theme", used by the notes master)
# so that theme1.xml ends up holding the default "Office Theme" colour
# scheme and theme2.xml ends up holding the former "Integral" colours.
#
# The only part of a theme that the PowerPoint object model exposes for
# in-place editing is the 12-slot colour scheme, reachable from any slide
# via Slide.ThemeColorScheme (this writes straight through to the theme
# part used by the presentation's slide master, i.e. ppt/theme/theme1.xml).
# We drive that API with the "Office Theme" RGB values so theme1.xml's
# colour scheme matches the post-edit target.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme colour scheme, in clrScheme slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
